# Update AgTests (F) and AgPosit (G) values for the Slovakia Covid daily
# stats sheet to reflect the "st 02. 06. 2021" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 423; F = 439035; G = 636 },
    @{ Row = 425; F = 138369; G = 548 },
    @{ Row = 426; F = 106673; G = 382 },
    @{ Row = 427; F = 90453;  G = 365 },
    @{ Row = 428; F = 102371 },
    @{ Row = 429; F = 178294; G = 458 },
    @{ Row = 430; F = 175386 },
    @{ Row = 432; F = 123675; G = 430 },
    @{ Row = 433; F = 86199;  G = 263 },
    @{ Row = 434; F = 79469;  G = 281 },
    @{ Row = 435; F = 82490 },
    @{ Row = 436; F = 144969; G = 353 },
    @{ Row = 437; F = 167197; G = 273 },
    @{ Row = 438; F = 121392; G = 250 },
    @{ Row = 439; F = 89095 },
    @{ Row = 440; F = 73595;  G = 226 },
    @{ Row = 442; F = 70109 },
    @{ Row = 443; F = 106691 },
    @{ Row = 444; F = 103351 },
    @{ Row = 449; F = 59412 },
    @{ Row = 450; F = 90553 },
    @{ Row = 451; F = 84942;  G = 114 },
    @{ Row = 452; F = 74176;  G = 124 },
    @{ Row = 453; F = 68881;  G = 209 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('F')) {
        $ws.Cells.Item($r, 6).Value = $u.F
    }
    if ($u.ContainsKey('G')) {
        $ws.Cells.Item($r, 7).Value = $u.G
    }
}

$wb.Save()
